{"js": "const replacements = [\n  [\n    \"2025-04-10 Thursday\",\n    \"2025-04-11 Friday\"\n  ],\n  [\n    \"46\u00f79=\",\n    \"88\u00f73=\"\n  ],\n  [\n    \"66\u00f75=\",\n    \"52\u00f77=\"\n  ],\n  [\n    \"69\u00f76=\",\n    \"87\u00f74=\"\n  ],\n  [\n    \"56\u00f72=\",\n    \"18\u00f78=\"\n  ],\n  [\n    \"87\u00f76=\",\n    \"68\u00f77=\"\n  ],\n  [\n    \"84\u00f76=\",\n    \"14\u00f79=\"\n  ],\n  [\n    \"62\u00f79=\",\n    \"43\u00f76=\"\n  ],\n  [\n    \"33\u00f79=\",\n    \"16\u00f78=\"\n  ],\n  [\n    \"68\u00f73=\",\n    \"67\u00f74=\"\n  ],\n  [\n    \"15\u00f78=\",\n    \"60\u00f76=\"\n  ],\n  [\n    \"47\u00f73=\",\n    \"98\u00f72=\"\n  ],\n  [\n    \"55\u00f76=\",\n    \"41\u00f72=\"\n  ],\n  [\n    \"20\u00f74=\",\n    \"49\u00f78=\"\n  ],\n  [\n    \"99\u00f74=\",\n    \"77\u00f72=\"\n  ],\n  [\n    \"19\u00f73=\",\n    \"63\u00f72=\"\n  ],\n  [\n    \"54\u00f78=\",\n    \"46\u00f73=\"\n  ],\n  [\n    \"32\u00f78=\",\n    \"44\u00f73=\"\n  ],\n  [\n    \"43\u00f79=\",\n    \"99\u00f76=\"\n  ],\n  [\n    \"28\u00f75=\",\n    \"94\u00f77=\"\n  ],\n  [\n    \"17\u00f79=\",\n    \"52\u00f73=\"\n  ],\n  [\n    \"91\u00f74=\",\n    \"63\u00f79=\"\n  ],\n  [\n    \"78\u00f72=\",\n    \"46\u00f72=\"\n  ],\n  [\n    \"53\u00f72=\",\n    \"60\u00f79=\"\n  ],\n  [\n    \"37\u00f75=\",\n    \"91\u00f75=\"\n  ],\n  [\n    \"33\u00f78=\",\n    \"63\u00f77=\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-04-10 Thursday', '2025-04-11 Friday'),\n    @('46\u00f79=', '88\u00f73='),\n    @('66\u00f75=', '52\u00f77='),\n    @('69\u00f76=', '87\u00f74='),\n    @('56\u00f72=', '18\u00f78='),\n    @('87\u00f76=', '68\u00f77='),\n    @('84\u00f76=', '14\u00f79='),\n    @('62\u00f79=', '43\u00f76='),\n    @('33\u00f79=', '16\u00f78='),\n    @('68\u00f73=', '67\u00f74='),\n    @('15\u00f78=', '60\u00f76='),\n    @('47\u00f73=', '98\u00f72='),\n    @('55\u00f76=', '41\u00f72='),\n    @('20\u00f74=', '49\u00f78='),\n    @('99\u00f74=', '77\u00f72='),\n    @('19\u00f73=', '63\u00f72='),\n    @('54\u00f78=', '46\u00f73='),\n    @('32\u00f78=', '44\u00f73='),\n    @('43\u00f79=', '99\u00f76='),\n    @('28\u00f75=', '94\u00f77='),\n    @('17\u00f79=', '52\u00f73='),\n    @('91\u00f74=', '63\u00f79='),\n    @('78\u00f72=', '46\u00f72='),\n    @('53\u00f72=', '60\u00f79='),\n    @('37\u00f75=', '91\u00f75='),\n    @('33\u00f78=', '63\u00f77='),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)\n}\n"}
